$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated D (Price) and E (Volume) values for rows where price changed.
$priceUpdates = @{
    2 = @("29.254.64", "  +0.53%  ")
    3 = @("1.857.66", "  +0.39%  ")
    4 = @("0.9996", "  +0.05%  ")
    5 = @("0.7068", "  +1.81%  ")
    6 = @("238.10", "  +0.11%  ")
    7 = @("1.000", "  +0.14%  ")
    8 = @("0.07997", "  +3.00%  ")
    9 = @("0.3027", "  -0.54%  ")
    10 = @("23.50", "  +0.81%  ")
    11 = @("0.08179", "  +0.78%  ")
    12 = @("1.865.85", "  +0.48%  ")
    13 = @("5.194", "  -0.33%  ")
    14 = @("0.7057", "  -2.78%  ")
    15 = @("89.70", "  +0.78%  ")
    16 = @("29.218.71", "  +0.41%  ")
    17 = @("0.000007933", "  +1.45%  ")
    18 = @("5.800", "  +0.90%  ")
    19 = @("13.24", "  +0.31%  ")
    20 = @("238.76", "  +1.21%  ")
    21 = @("0.9990", "  +0.01%  ")
    22 = @("2.093.01", "  -0.25%  ")
    23 = @("1.000", "  +0.14%  ")
    24 = @("7.475", "  -1.63%  ")
    25 = @("162.98", "  +1.20%  ")
    26 = @("8.875", "  -1.17%  ")
    27 = @("0.1433", "  -0.03%  ")
    28 = @("18.10", "  +0.23%  ")
    29 = @("1.926", "  -2.74%  ")
    31 = @("1.478", "  -0.60%  ")
    32 = @("4.372", "  -2.65%  ")
    33 = @("4.022", "  +0.39%  ")
    34 = @("0.05192", "  -0.77%  ")
    35 = @("1.160", "  -1.84%  ")
    36 = @("0.7139", "  +1.27%  ")
    38 = @("2.653", "  +0.29%  ")
    39 = @("0.01854", "  -0.10%  ")
    40 = @("2.723", "  +1.98%  ")
    41 = @("0.9365", "  +2.50%  ")
    42 = @("1.138.89", "  +4.04%  ")
    43 = @("5.948", "  -0.96%  ")
    44 = @("0.4265", "  -0.05%  ")
    45 = @("70.44", "  -0.33%  ")
    46 = @("0.9999", "  +0.14%  ")
    47 = @("102.65", "  -0.18%  ")
    48 = @("0.5304", "  -4.23%  ")
    49 = @("1.763", "  -0.30%  ")
    50 = @("1.999.54", "  +0.25%  ")
    51 = @("9.171", "  -0.02%  ")
}

# Rows where only the Volume(1h) value changed; Price stayed the same text.
$volumeOnlyUpdates = @{
    30 = "  +2.14%  "
    37 = "  -0.01%  "
}

foreach ($row in $priceUpdates.Keys) {
    $vals = $priceUpdates[$row]
    $dCell = $ws.Range("D$row")
    # Force text storage so numeric-looking strings (e.g. "1.000") are not
    # coerced into numbers, then restore the original General format so the
    # cell formatting matches the source workbook.
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $dCell.NumberFormat = "General"
    $ws.Range("E$row").Value = $vals[1]
}

foreach ($row in $volumeOnlyUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeOnlyUpdates[$row]
}
